$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, M, N, O, P, S for rows 2-20.
# (Weekly re-shuffle of dates/volumes/prices; other columns unaffected.)
$data = @{
    2  = @{ D = 44300; M = 45; N = 22000; O = 22000; P = 22000; S = 1100 }
    3  = @{ D = 44292; M = 30; N = 25000; O = 25000; P = 25000; S = 1250 }
    4  = @{ D = 44305; M = 20; N = 22000; O = 22000; P = 22000; S = 1100 }
    5  = @{ D = 44301; M = 38; N = 22000; O = 22000; P = 22000; S = 1100 }
    6  = @{ D = 44445; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    7  = @{ D = 44376; M = 38; N = 20000; O = 20000; P = 20000; S = 1000 }
    8  = @{ D = 44403; M = 50; N = 20000; O = 20000; P = 20000; S = 1000 }
    9  = @{ D = 44406; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    10 = @{ D = 44413; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    11 = @{ D = 44448; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    12 = @{ D = 44298; M = 65; N = 22000; O = 22000; P = 22000; S = 1100 }
    13 = @{ D = 44400; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    14 = @{ D = 44377; M = 25; N = 20000; O = 20000; P = 20000; S = 1000 }
    15 = @{ D = 44389; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    16 = @{ D = 44294; M = 25; N = 25000; O = 25000; P = 25000; S = 1250 }
    17 = @{ D = 44382; M = 24; N = 20000; O = 20000; P = 20000; S = 1000 }
    18 = @{ D = 44385; M = 36; N = 20000; O = 20000; P = 20000; S = 1000 }
    19 = @{ D = 44291; M = 70; N = 25000; O = 25000; P = 25000; S = 1250 }
    20 = @{ D = 44307; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
